# VAN-1811: Prepare and write FUNCTIONAL test cases and test scripts
#
# Updates the "Previous Doc" column (AX) for every data row (2-15) from the
# numeric placeholder 214002901789 to the text value 3703075487, and
# re-selects row 15 in the grid (mirrors clicking the row-15 header).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# --- Update AX2:AX15 ("Previous Doc") -> text "3703075487" -------------
$axRange = $ws.Range("AX2:AX15")
$axRange.Value = "'3703075487"

# Re-apply the original number format / style (column AX cells were style
# index 1, i.e. plain "Normal 2" cells) so only the value + type changes,
# keeping the cell's look-and-feel identical to before the edit.
$styleSource = $ws.Range("D2")
$styleSource.Copy() | Out-Null
$axRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Re-select entire row 15 (A15:XFD15), as if the row header was clicked
$ws.Rows(15).Select() | Out-Null
